# mise à jour aides pdf
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PDF links replacing the old .html/.htm support links, and refreshed
# "Date_lien" value (2025-02-08 = serial 45696) for every row.
$links = @(
    @{ Row = 2;  Url = "https://ductair.github.io/ductaironline/Support/Nouveautés 4x.pdf" },
    @{ Row = 3;  Url = "https://ductair.github.io/ductaironline/Support/MU.pdf" },
    @{ Row = 4;  Url = "https://ductair.github.io/ductaironline/Support/apropos.dpf" },
    @{ Row = 5;  Url = "https://ductair.github.io/ductaironline/Support/dournisseurs.pdf" },
    @{ Row = 6;  Url = "https://ductair.github.io/ductaironline/Support/Aide coudes.pdf" },
    @{ Row = 7;  Url = "https://ductair.github.io/ductaironline/Support/Aide_projection.pdf" },
    @{ Row = 8;  Url = "https://ductair.github.io/ductaironline/Support/menu.pdf" },
    @{ Row = 9;  Url = "https://ductair.github.io/ductaironline/Remote/Commun/ductair.png" }
)

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 3)   # column C = "Lien"
    $cell.Value2 = $link.Url
    $ws.Hyperlinks.Add($cell, $link.Url) | Out-Null
    $ws.Cells.Item($link.Row, 4).Value = 45696   # column D = "Date_lien"
}

# Column B ("Aide") grew a bit wider once the sheet was touched again.
$ws.Columns(2).ColumnWidth = 16.5

# Leave the selection where the author ended up after the edits.
$ws.Range("C10").Select() | Out-Null
